$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.943.81"
$ws.Range("E2").Value = "  +0.58%  "

$ws.Range("D3").Value = "1.587.41"
$ws.Range("E3").Value = "  +0.14%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "'210.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.50%  "

$ws.Range("E6").Value = "  -0.14%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "'0.246"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.44%  "

$ws.Range("E9").Value = "  -1.13%  "

$ws.Range("D10").Value = "'17.81"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.17%  "

$ws.Range("D11").Value = "'0.0809"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.47%  "

$ws.Range("D12").Value = "1.811.70"
$ws.Range("E12").Value = "  +0.37%  "

$ws.Range("D13").Value = "1.589.62"
$ws.Range("E13").Value = "  +0.27%  "

$ws.Range("D14").Value = "'3.97"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.28%  "

$ws.Range("D15").Value = "'0.510"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.04%  "

$ws.Range("D16").Value = "25.935.16"
$ws.Range("E16").Value = "  +0.57%  "

$ws.Range("D17").Value = "'59.88"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.12%  "

$ws.Range("D18").Value = "0.0₃0718"
$ws.Range("E18").Value = "  -0.38%  "

$ws.Range("E19").Value = "  -0.22%  "

$ws.Range("D20").Value = "'198.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.34%  "

$ws.Range("E21").Value = "  +0.40%  "

$ws.Range("D22").Value = "'9.18"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.79%  "

$ws.Range("D23").Value = "'5.95"
$ws.Range("D23").Style = "Normal"

$ws.Range("E24").Value = "  +8.72%  "

$ws.Range("D25").Value = "'143.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.10%  "

$ws.Range("E26").Value = "  -0.15%  "

$ws.Range("E27").Value = "  -8.28%  "

$ws.Range("D28").Value = "'15.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.32%  "

$ws.Range("E29").Value = "  -0.23%  "

$ws.Range("E31").Value = "  +0.41%  "

$ws.Range("E32").Value = "  +0.07%  "

$ws.Range("D33").Value = "'2.93"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.70%  "

$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "'2.37"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.56%  "

$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "'1.46"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.74%  "

$ws.Range("D36").Value = "1.121.97"
$ws.Range("E36").Value = "  +2.20%  "

$ws.Range("E37").Value = "  +7.86%  "

$ws.Range("E38").Value = "  -0.15%  "

$ws.Range("E39").Value = "  -0.88%  "

$ws.Range("D40").Value = "'0.777"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.23%  "

$ws.Range("D41").Value = "'0.487"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.41%  "

$ws.Range("D42").Value = "'0.781"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.42%  "

$ws.Range("D43").Value = "1.723.43"
$ws.Range("E43").Value = "  +0.22%  "

$ws.Range("D44").Value = "'5.07"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.51%  "

$ws.Range("D45").Value = "'91.71"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.15%  "

$ws.Range("D46").Value = "'1.48"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.19%  "

$ws.Range("D47").Value = "'53.07"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.11%  "

$ws.Range("E48").Value = "  -1.21%  "

$ws.Range("E49").Value = "  -0.25%  "

$ws.Range("E50").Value = "  +0.23%  "

$ws.Range("D51").Value = "0.0₇0918"
$ws.Range("E51").Value = "  -17.88%  "
